$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "Summer Day"
$ws.Range("F4").Value = "Summer Day Kevin MacLeod (incompetech.com)"
$ws.Range("F5").Value = "Licensed under Creative Commons: By Attribution 3.0 License"
$ws.Range("F6").Value = "http://creativecommons.org/licenses/by/3.0/"

$ws.Range("B7").Value = "Kevin Macleod"
$ws.Range("C7").Value = "https://www.youtube.com/audiolibrary/music"

$noteRange = $ws.Range("F4:F6")
$noteRange.Font.Bold = $true
$noteRange.Font.Size = 8
$noteRange.Font.Name = "Arial"
$noteRange.Font.Color = 0

$ws.Range("C7").Select()
